$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet keeps the last populated row styled differently (plain date
# format) than the rows above it (date + time format). Appending a new
# row means the old last row (11) reverts to the "normal" style, and the
# newly appended row (12) becomes the new "last row" with the special
# style that row 11 used to have.

# 1) Copy row 11's current ("last row") formatting onto the new row 12
#    before row 11's own formatting changes.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

# 2) Row 11 reverts to the same formatting used by the rows above it.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

# 3) Fill in the new row's values.
$ws.Range("A12").Value = 45752
$ws.Range("B12").Value = 42
$ws.Range("C12").Value = 43
$ws.Range("D12").Value = 43
